# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-market-day records (rows 2-19) of the
# "Hortaliza, Vega Monumental Concepción - Espárragos" sheet: the data that
# used to live on one row now lives on another row (same set of records,
# new row order/assignment). Columns A,B,C,E,F,G,Q,R are identical on every
# row already, so only D,H,I,J,K,L,M,N,O,P need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 19

# Columns whose values vary row to row and therefore need to be relocated.
$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "O", "P")

# Snapshot every relevant cell's current value before writing anything,
# so that writes to one row never clobber data still needed for another.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# Destination row -> source row (which row's original data now appears there).
$sourceForRow = @{
    2  = 17
    3  = 18
    4  = 6
    5  = 3
    6  = 4
    7  = 19
    8  = 7
    9  = 12
    10 = 14
    11 = 9
    12 = 11
    13 = 8
    14 = 5
    15 = 2
    16 = 13
    17 = 15
    18 = 16
    19 = 10
}

for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $sourceForRow[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
